# Update need_to_buy.xlsx forecast data (columns B, C, F) with refreshed
# values produced by the upstream R script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> @(B, C, F) new values ($null means "leave unchanged")
$updates = @{
    2  = @(7113.78710552895, 6829.17745684863, 2.14099534665335)
    3  = @(7122.41926747084, 6771.12523143195, 290.124665415046)
    4  = @(6720.65507183883, 6393.30885483788, 266.795167416627)
    5  = @(6432.57135429943, 6130.18866320353, 246.883818829338)
    6  = @(6366.54402734252, 5645.04376774423, 211.383259933405)
    7  = @(2137.91168797406, 4135.71432723904, 147.316402511041)
    9  = @($null,            7102.78408619527, 310.119694872751)
    10 = @($null,            7593.36854654761, 330.560714054098)
    11 = @($null,            7302.78918035739, 318.45324046284)
    12 = @($null,            7219.66881486471, 314.989891900644)
    13 = @($null,            5846.72470409775, 253.522731679543)
    14 = @($null,            4234.03298007117, 60.0611005049032)
    15 = @($null,            4213.59292188353, 59.2094314137514)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($null -ne $vals[0]) {
        $ws.Cells.Item($row, 2).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($row, 3).Value = $vals[1]
    }
    if ($null -ne $vals[2]) {
        $ws.Cells.Item($row, 6).Value = $vals[2]
    }
}
